# Auto-generated edit script updating cryptos price/volume cells
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.938.10"
$ws.Range("E2").Value = "  -0.70%  "
$ws.Range("D3").Value = "1.645.42"
$ws.Range("E3").Value = "  -0.90%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "'310.15"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("E7").Value = "  -0.72%  "
$ws.Range("D8").Value = "'0.3830"
$ws.Range("E8").Value = "  -1.42%  "
$ws.Range("D9").Value = "'50.66"
$ws.Range("E9").Value = "  -1.90%  "
$ws.Range("D10").Value = "'1.333"
$ws.Range("E10").Value = "  -3.34%  "
$ws.Range("E11").Value = "  +0.43%  "
$ws.Range("D12").Value = "'0.08405"
$ws.Range("E12").Value = "  -1.10%  "
$ws.Range("D13").Value = "'23.81"
$ws.Range("E13").Value = "  -1.83%  "
$ws.Range("D14").Value = "'6.989"
$ws.Range("E14").Value = "  -3.96%  "
$ws.Range("D15").Value = "'7.846"
$ws.Range("E15").Value = "  -3.77%  "
$ws.Range("D16").Value = "'0.00001313"
$ws.Range("E16").Value = "  -0.66%  "
$ws.Range("D17").Value = "1.645.60"
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("D18").Value = "'93.92"
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("D19").Value = "'0.06957"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").Value = "'19.50"
$ws.Range("E20").Value = "  -3.06%  "
$ws.Range("D21").Value = "'6.896"
$ws.Range("E21").Value = "  -1.17%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("E23").Value = "  -1.02%  "
$ws.Range("D24").Value = "23.908.35"
$ws.Range("E24").Value = "  -0.76%  "
$ws.Range("D25").Value = "'2.435"
$ws.Range("E25").Value = "  -2.30%  "
$ws.Range("D26").Value = "'2.899"
$ws.Range("E26").Value = "  -8.68%  "
$ws.Range("D27").Value = "'21.93"
$ws.Range("E27").Value = "  -1.53%  "
$ws.Range("D28").Value = "'153.31"
$ws.Range("E28").Value = "  -0.42%  "
$ws.Range("D29").Value = "'5.576"
$ws.Range("E29").Value = "  +5.19%  "
$ws.Range("D30").Value = "'136.97"
$ws.Range("E30").Value = "  -1.97%  "
$ws.Range("D31").Value = "'2.501"
$ws.Range("E31").Value = "  +0.41%  "
$ws.Range("D32").Value = "'7.633"
$ws.Range("E32").Value = "  -2.98%  "
$ws.Range("D33").Value = "1.828.18"
$ws.Range("E33").Value = "  -0.36%  "
$ws.Range("D34").Value = "'0.08042"
$ws.Range("E34").Value = "  -1.57%  "
$ws.Range("D35").Value = "'0.9810"
$ws.Range("E35").Value = "  -7.01%  "
$ws.Range("D36").Value = "'0.02910"
$ws.Range("E36").Value = "  -4.11%  "
$ws.Range("D37").Value = "'6.597"
$ws.Range("E37").Value = "  -2.39%  "
$ws.Range("D38").Value = "'0.2670"
$ws.Range("E38").Value = "  -2.67%  "
$ws.Range("D39").Value = "'10.39"
$ws.Range("E39").Value = "  -8.02%  "
$ws.Range("D40").Value = "'0.09096"
$ws.Range("E40").Value = "  -1.14%  "
$ws.Range("D41").Value = "'0.7516"
$ws.Range("E41").Value = "  -1.27%  "
$ws.Range("D42").Value = "'13.34"
$ws.Range("E42").Value = "  -1.60%  "
$ws.Range("D43").Value = "'1.418"
$ws.Range("E43").Value = "  -0.80%  "
$ws.Range("D44").Value = "'16.51"
$ws.Range("E44").Value = "  -1.11%  "
$ws.Range("D45").Value = "'0.6910"
$ws.Range("D46").Value = "'2.422"
$ws.Range("E46").Value = "  -3.54%  "
$ws.Range("D47").Value = "'4.090"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("D49").Value = "'0.08261"
$ws.Range("E49").Value = "  -1.41%  "
$ws.Range("D50").Value = "'134.69"
$ws.Range("E50").Value = "  -0.85%  "
$ws.Range("D51").Value = "'1.222"
$ws.Range("E51").Value = "  -1.97%  "
